$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.125.78"
$ws.Range("E2").Value = "  -1.34%  "
$ws.Range("D3").Value = "1.835.40"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6647"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.55%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2952"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07349"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.61%  "
$ws.Range("E10").Value = "  -4.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07678"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.51%  "
$ws.Range("D12").Value = "1.845.52"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.020"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6756"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "86.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.204"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.47%  "
$ws.Range("D17").Value = "29.059.21"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008226"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "229.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.02%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.314"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "160.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1418"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.694"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.499"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.225"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.094"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.194"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05299"
$ws.Range("D32").Style = "Normal"
$ws.Range("E33").Value = "  -1.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7481"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.129"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.681"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").Value = "1.317.07"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01804"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.717"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9210"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.981"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9983"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "103.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.30%  "
$ws.Range("D44").Value = "1.986.01"
$ws.Range("E44").Value = "  -0.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5164"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.30%  "
$ws.Range("E46").Value = "  -3.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "63.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.760"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.282"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05934"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07304"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.87%  "
